# Creation of html pages
# The only durable content change in this commit is the worksheet
# tab name: "Feuil1" -> "ST4" (everything else in the diff is Excel
# session/save metadata - file version build, absolute source path,
# revision GUIDs, window geometry - that Excel regenerates on every
# save and isn't a deliberate, scriptable edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Name = "ST4"
